$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08758366666666667
$ws.Range("H2").Value = 0.262751
$ws.Range("I2").Value = 0.2371976925785164
$ws.Range("J2").Value = 0.2371976925785164
$ws.Range("M2").Value = 176.8550973333333
$ws.Range("N2").Value = 530.565292
$ws.Range("O2").Value = 0.2669710696905332
$ws.Range("P2").Value = 0.2669710696905332
$ws.Range("Q2").Value = 15.48961789314356
$ws.Range("R2").Value = 139.406561038292
$ws.Range("S2").Value = 0.06332492171581278
$ws.Range("T2").Value = 0.06332492171581278

$ws.Range("G3").Value = 0.08758366666666667
$ws.Range("H3").Value = 0.262751
$ws.Range("I3").Value = 0.2371976925785164
$ws.Range("J3").Value = 0.2371976925785164
$ws.Range("O3").Value = 0.2198524722701247
$ws.Range("P3").Value = 0.2198524722701247
$ws.Range("Q3").Value = 12.75580456067644
$ws.Range("R3").Value = 114.802241046088
$ws.Range("S3").Value = 0.05214849913015585
$ws.Range("T3").Value = 0.05214849913015585

$ws.Range("G4").Value = 0.08758366666666667
$ws.Range("H4").Value = 0.262751
$ws.Range("I4").Value = 0.2371976925785164
$ws.Range("J4").Value = 0.2371976925785164
$ws.Range("M4").Value = 84.02511333333334
$ws.Range("N4").Value = 252.07534
$ws.Range("O4").Value = 0.1268398520919549
$ws.Range("P4").Value = 0.1268398520919549
$ws.Range("Q4").Value = 7.359227517815556
$ws.Range("R4").Value = 66.23304766034001
$ws.Range("S4").Value = 0.030086120243212
$ws.Range("T4").Value = 0.030086120243212

$ws.Range("G5").Value = 0.08758366666666667
$ws.Range("H5").Value = 0.262751
$ws.Range("I5").Value = 0.2371976925785164
$ws.Range("J5").Value = 0.2371976925785164
$ws.Range("M5").Value = 84.92877566666668
$ws.Range("N5").Value = 254.786327
$ws.Range("O5").Value = 0.1282039727953256
$ws.Range("P5").Value = 0.1282039727953256
$ws.Range("Q5").Value = 7.438373578397446
$ws.Range("R5").Value = 66.94536220557701
$ws.Range("S5").Value = 0.03040968652645012
$ws.Range("T5").Value = 0.03040968652645012

$ws.Range("G6").Value = 0.08758366666666667
$ws.Range("H6").Value = 0.262751
$ws.Range("I6").Value = 0.2371976925785164
$ws.Range("J6").Value = 0.2371976925785164
$ws.Range("M6").Value = 63.97102366666667
$ws.Range("N6").Value = 191.913071
$ws.Range("O6").Value = 0.09656726254996952
$ws.Range("P6").Value = 0.09656726254996952
$ws.Range("Q6").Value = 5.602816813146778
$ws.Range("R6").Value = 50.42535131832101
$ws.Range("S6").Value = 0.02290553185547655
$ws.Range("T6").Value = 0.02290553185547655

$ws.Range("G7").Value = 0.08758366666666667
$ws.Range("H7").Value = 0.262751
$ws.Range("I7").Value = 0.2371976925785164
$ws.Range("J7").Value = 0.2371976925785164
$ws.Range("M7").Value = 107.0290476666667
$ws.Range("N7").Value = 321.087143
$ws.Range("O7").Value = 0.1615653706020921
$ws.Range("P7").Value = 0.1615653706020921
$ws.Range("Q7").Value = 9.373996434488111
$ws.Range("R7").Value = 84.36596791039301
$ws.Range("S7").Value = 0.03832293310740911
$ws.Range("T7").Value = 0.03832293310740911

$ws.Range("I8").Value = 0.4417069141397272
$ws.Range("J8").Value = 0.4417069141397272
$ws.Range("M8").Value = 176.8550973333333
$ws.Range("N8").Value = 530.565292
$ws.Range("O8").Value = 0.2669710696905332
$ws.Range("P8").Value = 0.2669710696905332
$ws.Range("Q8").Value = 28.84459476147378
$ws.Range("R8").Value = 259.601352853264
$ws.Range("S8").Value = 0.1179229673575875
$ws.Range("T8").Value = 0.1179229673575875

$ws.Range("I9").Value = 0.4417069141397272
$ws.Range("J9").Value = 0.4417069141397272
$ws.Range("O9").Value = 0.2198524722701247
$ws.Range("P9").Value = 0.2198524722701247
$ws.Range("S9").Value = 0.09711035709242674
$ws.Range("T9").Value = 0.09711035709242674

$ws.Range("I10").Value = 0.4417069141397272
$ws.Range("J10").Value = 0.4417069141397272
$ws.Range("M10").Value = 84.02511333333334
$ws.Range("N10").Value = 252.07534
$ws.Range("O10").Value = 0.1268398520919549
$ws.Range("P10").Value = 0.1268398520919549
$ws.Range("Q10").Value = 13.70427191769778
$ws.Range("R10").Value = 123.33844725928
$ws.Range("S10").Value = 0.0560260396574768
$ws.Range("T10").Value = 0.0560260396574768

$ws.Range("I11").Value = 0.4417069141397272
$ws.Range("J11").Value = 0.4417069141397272
$ws.Range("M11").Value = 84.92877566666668
$ws.Range("N11").Value = 254.786327
$ws.Range("O11").Value = 0.1282039727953256
$ws.Range("P11").Value = 0.1282039727953256
$ws.Range("Q11").Value = 13.85165683449823
$ws.Range("R11").Value = 124.664911510484
$ws.Range("S11").Value = 0.05662858120387681
$ws.Range("T11").Value = 0.05662858120387681

$ws.Range("I12").Value = 0.4417069141397272
$ws.Range("J12").Value = 0.4417069141397272
$ws.Range("M12").Value = 63.97102366666667
$ws.Range("N12").Value = 191.913071
$ws.Range("O12").Value = 0.09656726254996952
$ws.Range("P12").Value = 0.09656726254996952
$ws.Range("Q12").Value = 10.43350337063689
$ws.Range("R12").Value = 93.90153033573201
$ws.Range("S12").Value = 0.04265442754786788
$ws.Range("T12").Value = 0.04265442754786788

$ws.Range("I13").Value = 0.4417069141397272
$ws.Range("J13").Value = 0.4417069141397272
$ws.Range("M13").Value = 107.0290476666667
$ws.Range("N13").Value = 321.087143
$ws.Range("O13").Value = 0.1615653706020921
$ws.Range("P13").Value = 0.1615653706020921
$ws.Range("Q13").Value = 17.45615226363956
$ws.Range("R13").Value = 157.105370372756
$ws.Range("S13").Value = 0.07136454128049149
$ws.Range("T13").Value = 0.07136454128049149

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1185623333333333
$ws.Range("H14").Value = 0.355687
$ws.Range("I14").Value = 0.3210953932817563
$ws.Range("J14").Value = 0.3210953932817564
$ws.Range("M14").Value = 176.8550973333333
$ws.Range("N14").Value = 530.565292
$ws.Range("O14").Value = 0.2669710696905332
$ws.Range("P14").Value = 0.2669710696905332
$ws.Range("Q14").Value = 20.96835300173378
$ws.Range("R14").Value = 188.715177015604
$ws.Range("S14").Value = 0.08572318061713294
$ws.Range("T14").Value = 0.08572318061713295

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1185623333333333
$ws.Range("H15").Value = 0.355687
$ws.Range("I15").Value = 0.3210953932817563
$ws.Range("J15").Value = 0.3210953932817564
$ws.Range("O15").Value = 0.2198524722701247
$ws.Range("P15").Value = 0.2198524722701247
$ws.Range("Q15").Value = 17.26757978760622
$ws.Range("R15").Value = 155.408218088456
$ws.Range("S15").Value = 0.07059361604754213
$ws.Range("T15").Value = 0.07059361604754215

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1185623333333333
$ws.Range("H16").Value = 0.355687
$ws.Range("I16").Value = 0.3210953932817563
$ws.Range("J16").Value = 0.3210953932817564
$ws.Range("M16").Value = 84.02511333333334
$ws.Range("N16").Value = 252.07534
$ws.Range("O16").Value = 0.1268398520919549
$ws.Range("P16").Value = 0.1268398520919549
$ws.Range("Q16").Value = 9.962213495397778
$ws.Range("R16").Value = 89.65992145858002
$ws.Range("S16").Value = 0.04072769219126605
$ws.Range("T16").Value = 0.04072769219126606

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1185623333333333
$ws.Range("H17").Value = 0.355687
$ws.Range("I17").Value = 0.3210953932817563
$ws.Range("J17").Value = 0.3210953932817564
$ws.Range("M17").Value = 84.92877566666668
$ws.Range("N17").Value = 254.786327
$ws.Range("O17").Value = 0.1282039727953256
$ws.Range("P17").Value = 0.1282039727953256
$ws.Range("Q17").Value = 10.06935381018322
$ws.Range("R17").Value = 90.62418429164902
$ws.Range("S17").Value = 0.04116570506499867
$ws.Range("T17").Value = 0.04116570506499868

$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.1185623333333333
$ws.Range("H18").Value = 0.355687
$ws.Range("I18").Value = 0.3210953932817563
$ws.Range("J18").Value = 0.3210953932817564
$ws.Range("M18").Value = 63.97102366666667
$ws.Range("N18").Value = 191.913071
$ws.Range("O18").Value = 0.09656726254996952
$ws.Range("P18").Value = 0.09656726254996952
$ws.Range("Q18").Value = 7.584553831641889
$ws.Range("R18").Value = 68.26098448477701
$ws.Range("S18").Value = 0.03100730314662508
$ws.Range("T18").Value = 0.03100730314662509

$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.1185623333333333
$ws.Range("H19").Value = 0.355687
$ws.Range("I19").Value = 0.3210953932817563
$ws.Range("J19").Value = 0.3210953932817564
$ws.Range("M19").Value = 107.0290476666667
$ws.Range("N19").Value = 321.087143
$ws.Range("O19").Value = 0.1615653706020921
$ws.Range("P19").Value = 0.1615653706020921
$ws.Range("Q19").Value = 12.68961362580456
$ws.Range("R19").Value = 114.206522632241
$ws.Range("S19").Value = 0.05187789621419147
$ws.Range("T19").Value = 0.05187789621419148
